$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column E header + data (LWW JAGS model results)
$ws.Range("E1").Value = "JAGS LWW"
$ws.Range("E1").Font.Bold = $true

$ws.Range("E15").Value = "3.5 (1.5-7.2)"
$ws.Range("E14").Value = "10.0 (7.8-12.9)"
$ws.Range("E16").Value = "4.3(1.3-7.0)"
$ws.Range("E17").Value = "9.0 (6.5-11.4)"
$ws.Range("E18").Value = "13.8 (10.6-21.9)"
$ws.Range("E8").Value = "10.9 (7.2-18.6)"
$ws.Range("E9").Value = "2.0 (0.6-4.5)"
$ws.Range("E10").Value = "2.4 (0.1-5.9)"
$ws.Range("E11").Value = "9.1 (4.8-13.3)"
$ws.Range("E12").Value = "19.0 (12.7-89.8)"
$ws.Range("E3").Value = "5.0 (4.4, 7.7)"
$ws.Range("E4").Value = "1.4 (1.22,2.0)"
$ws.Range("E5").Value = "3.2 (1.6, 4.6)"
$ws.Range("E6").Value = "10.7 (7.6-19.8)"

# Column widths adjusted to fit new content
$ws.Columns.Item(3).ColumnWidth = 14.7109375
$ws.Columns.Item(4).ColumnWidth = 26.28515625
$ws.Columns.Item(5).ColumnWidth = 14.42578125

# Row 27's lone formatted cell moved down to rows 28/29/32
$ws.Range("D27").Clear()

# Extra formatted (blank) cells added alongside the new column / extended table
$ws.Range("K14").NumberFormat = "0.00%"
$ws.Range("L14").NumberFormat = "0.00%"
$ws.Range("I17").NumberFormat = "0%"
$ws.Range("I18").NumberFormat = "0%"
$ws.Range("I21").NumberFormat = "0%"
$ws.Range("I22").NumberFormat = "0%"
$ws.Range("F25").NumberFormat = "0.00%"
$ws.Range("G25").NumberFormat = "0.00%"
$ws.Range("D28").NumberFormat = "0%"
$ws.Range("D29").NumberFormat = "0%"
$ws.Range("D32").NumberFormat = "0%"

# Selection moved to N13 to match the post-edit cursor position
$ws.Range("N13").Select()
